$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,14
$data[0,0] = 1.113143423883969
$data[0,1] = 0.2653531116876309
$data[0,2] = 0
$data[0,3] = 0.2269633380091669
$data[0,4] = 1.82366755224642
$data[0,5] = 0.002449352365023655
$data[0,6] = 0
$data[0,7] = 0.802839715008556
$data[0,8] = 0.05239605534474379
$data[0,9] = 0
$data[0,10] = 0.4833106585977447
$data[0,11] = 0.330598998055784
$data[0,12] = 0
$data[0,13] = 2.963893769994129
$data[1,0] = 1.00875966013723
$data[1,1] = 0.2526534927238515
$data[1,2] = 0
$data[1,3] = 0.2283157693032791
$data[1,4] = 1.83257478891781
$data[1,5] = 0.002452060249102687
$data[1,6] = 0
$data[1,7] = 0.8146725878548011
$data[1,8] = 0.05026012193944496
$data[1,9] = 0
$data[1,10] = 0.4770610312162376
$data[1,11] = 0.3116502185809864
$data[1,12] = 0
$data[1,13] = 2.99206698316695
$data[2,0] = 0.9446429269162877
$data[2,1] = 0.24480957004441
$data[2,2] = 0
$data[2,3] = 0.2292089849096932
$data[2,4] = 1.839090532117787
$data[2,5] = 0.002453811730133503
$data[2,6] = 0
$data[2,7] = 0.8224513195235126
$data[2,8] = 0.04894679596576879
$data[2,9] = 0
$data[2,10] = 0.4733746686494698
$data[2,11] = 0.300068353366477
$data[2,12] = 0
$data[2,13] = 3.011454396790171
$data[3,0] = 0.9185103843911975
$data[3,1] = 0.2416016587374656
$data[3,2] = 0
$data[3,3] = 0.2295888027158037
$data[3,4] = 1.842008979656306
$data[3,5] = 0.002454547875825606
$data[3,6] = 0
$data[3,7] = 0.8257501801982237
$data[3,8] = 0.048411180936764
$data[3,9] = 0
$data[3,10] = 0.471910590635062
$data[3,11] = 0.2953622613861313
$data[3,12] = 0
$data[3,13] = 3.019879653926509
$data[4,0] = 0.9141708785122091
$data[4,1] = 0.2410683017782418
$data[4,2] = 0
$data[4,3] = 0.2296528278361505
$data[4,4] = 1.842509485139963
$data[4,5] = 0.002454671467126072
$data[4,6] = 0
$data[4,7] = 0.8263057394293298
$data[4,8] = 0.04832221810688964
$data[4,9] = 0
$data[4,10] = 0.4716697910760388
$data[4,11] = 0.2945816507261299
$data[4,12] = 0
$data[4,13] = 3.021310337517093
$data[5,0] = 0.9442905100059988
$data[5,1] = 0.2447663531268063
$data[5,2] = 0
$data[5,3] = 0.2292140431507104
$data[5,4] = 1.839128825468933
$data[5,5] = 0.002453821567092257
$data[5,6] = 0
$data[5,7] = 0.8224952870864826
$data[5,8] = 0.04893957413517924
$data[5,9] = 0
$data[5,10] = 0.4733547688803839
$data[5,11] = 0.3000048297904812
$data[5,12] = 0
$data[5,13] = 3.011565898889202
$data[6,0] = 1.077158039850019
$data[6,1] = 0.2609840062850992
$data[6,2] = 0
$data[6,3] = 0.2274166361918013
$data[6,4] = 1.826521552243399
$data[6,5] = 0.002450267646907911
$data[6,6] = 0
$data[6,7] = 0.8068130896847414
$data[6,8] = 0.05165999493734574
$data[6,9] = 0
$data[6,10] = 0.4811245724343394
$data[6,11] = 0.324054708116158
$data[6,12] = 0
$data[6,13] = 2.973173999076252
$data[7,0] = 1.337451232757473
$data[7,1] = 0.2924123275631132
$data[7,2] = 0
$data[7,3] = 0.224389049115997
$data[7,4] = 1.810103303341535
$data[7,5] = 0.002444000210704052
$data[7,6] = 0
$data[7,7] = 0.7801377653054544
$data[7,8] = 0.05697843096520927
$data[7,9] = 0
$data[7,10] = 0.49755146535351
$data[7,11] = 0.3716224634878174
$data[7,12] = 0
$data[7,13] = 2.914488497289568
$data[8,0] = 1.528464062253931
$data[8,1] = 0.3152672224410367
$data[8,2] = 0
$data[8,3] = 0.2224659577879695
$data[8,4] = 1.803105278263743
$data[8,5] = 0.002439819105266188
$data[8,6] = 0
$data[8,7] = 0.7630310530553643
$data[8,8] = 0.06087420600434257
$data[8,9] = 0
$data[8,10] = 0.510337874178262
$data[8,11] = 0.4068047518327376
$data[8,12] = 0
$data[8,13] = 2.881530186508456
$data[9,0] = 1.615299300942809
$data[9,1] = 0.3256120033708498
$data[9,2] = 0
$data[9,3] = 0.2216561361554881
$data[9,4] = 1.801021915756124
$data[9,5] = 0.002438008083947485
$data[9,6] = 0
$data[9,7] = 0.7557910712695097
$data[9,8] = 0.06264361751204461
$data[9,9] = 0
$data[9,10] = 0.516308956145906
$data[9,11] = 0.4228583997977182
$data[9,12] = 0
$data[9,13] = 2.868750480347074
$data[10,0] = 1.648171813406009
$data[10,1] = 0.3295216602717801
$data[10,2] = 0
$data[10,3] = 0.2213587956741101
$data[10,4] = 1.800391198802529
$data[10,5] = 0.002437335312475567
$data[10,6] = 0
$data[10,7] = 0.7531275304012759
$data[10,8] = 0.06331320883512603
$data[10,9] = 0
$data[10,10] = 0.5185920958910657
$data[10,11] = 0.4289442597687128
$data[10,12] = 0
$data[10,13] = 2.864230044061912
$data[11,0] = 1.64109260847016
$data[11,1] = 0.3286799908496505
$data[11,2] = 0
$data[11,3] = 0.2214224191059042
$data[11,4] = 1.800519998407296
$data[11,5] = 0.00243747962770508
$data[11,6] = 0
$data[11,7] = 0.7536976964034352
$data[11,8] = 0.06316902091988652
$data[11,9] = 0
$data[11,10] = 0.5180994045981322
$data[11,11] = 0.4276332690367255
$data[11,12] = 0
$data[11,13] = 2.865189403804095
$data[12,0] = 1.618003955425081
$data[12,1] = 0.3259338084307899
$data[12,2] = 0
$data[12,3] = 0.2216314871081444
$data[12,4] = 1.800966855487758
$data[12,5] = 0.002437952474103846
$data[12,6] = 0
$data[12,7] = 0.7555703743893076
$data[12,8] = 0.06269871439277352
$data[12,9] = 0
$data[12,10] = 0.5164963512072944
$data[12,11] = 0.4233589551054493
$data[12,12] = 0
$data[12,13] = 2.868372184843963
$data[13,0] = 1.603860117885063
$data[13,1] = 0.3242506852358815
$data[13,2] = 0
$data[13,3] = 0.221760760451696
$data[13,4] = 1.801261172105043
$data[13,5] = 0.002438243800347239
$data[13,6] = 0
$data[13,7] = 0.7567276169521726
$data[13,8] = 0.06241057861817723
$data[13,9] = 0
$data[13,10] = 0.5155172967193522
$data[13,11] = 0.4207416752060453
$data[13,12] = 0
$data[13,13] = 2.870363289187026
$data[14,0] = 1.522787860922278
$data[14,1] = 0.3145901021437112
$data[14,2] = 0
$data[14,3] = 0.2225201873356539
$data[14,4] = 1.803263568773417
$data[14,5] = 0.002439939285456908
$data[14,6] = 0
$data[14,7] = 0.7635151212351197
$data[14,8] = 0.06075851073887151
$data[14,9] = 0
$data[14,10] = 0.5099507419118368
$data[14,11] = 0.4057565642331227
$data[14,12] = 0
$data[14,13] = 2.882409965956782
$data[15,0] = 1.473036623588087
$data[15,1] = 0.3086501752074753
$data[15,2] = 0
$data[15,3] = 0.2230027013422866
$data[15,4] = 1.804773742893559
$data[15,5] = 0.002441002670717363
$data[15,6] = 0
$data[15,7] = 0.7678179490365942
$data[15,8] = 0.05974427280484207
$data[15,9] = 0
$data[15,10] = 0.5065752690599368
$data[15,11] = 0.396575988758407
$data[15,12] = 0
$data[15,13] = 2.890367567901876
$data[16,0] = 1.444415729179241
$data[16,1] = 0.3052288016466775
$data[16,2] = 0
$data[16,3] = 0.2232863504965188
$data[16,4] = 1.80574589832284
$data[16,5] = 0.002441622869373832
$data[16,6] = 0
$data[16,7] = 0.770343823378667
$data[16,8] = 0.05916064941778387
$data[16,9] = 0
$data[16,10] = 0.5046483383071063
$data[16,11] = 0.3913002026850165
$data[16,12] = 0
$data[16,13] = 2.895152844398609
$data[17,0] = 1.434724332553515
$data[17,1] = 0.304069550703872
$data[17,2] = 0
$data[17,3] = 0.2233834411074529
$data[17,4] = 1.806092836527057
$data[17,5] = 0.002441834331426195
$data[17,6] = 0
$data[17,7] = 0.7712077959791763
$data[17,8] = 0.0589630011729696
$data[17,9] = 0
$data[17,10] = 0.5039984178834516
$data[17,11] = 0.389514719818429
$data[17,12] = 0
$data[17,13] = 2.896808806546517
$data[18,0] = 1.478333289531179
$data[18,1] = 0.309282997615469
$data[18,2] = 0
$data[18,3] = 0.2229507037077703
$data[18,4] = 1.804602265809933
$data[18,5] = 0.002440888585477498
$data[18,6] = 0
$data[18,7] = 0.7673546260963242
$data[18,8] = 0.05985226745269046
$data[18,9] = 0
$data[18,10] = 0.5069330892398227
$data[18,11] = 0.3975527992861245
$data[18,12] = 0
$data[18,13] = 2.889498906844864
$data[19,0] = 1.624785941465916
$data[19,1] = 0.3267406390342273
$data[19,2] = 0
$data[19,3] = 0.2215698259935532
$data[19,4] = 1.80083130889831
$data[19,5] = 0.002437813234765628
$data[19,6] = 0
$data[19,7] = 0.7550182031093371
$data[19,8] = 0.06283686728733784
$data[19,9] = 0
$data[19,10] = 0.5169666107524051
$data[19,11] = 0.424614245968904
$data[19,12] = 0
$data[19,13] = 2.867428662775609
$data[20,0] = 1.720441452784826
$data[20,1] = 0.3381053162801777
$data[20,2] = 0
$data[20,3] = 0.2207216628858379
$data[20,4] = 1.799288930466929
$data[20,5] = 0.002435879195652992
$data[20,6] = 0
$data[20,7] = 0.7474108358853861
$data[20,8] = 0.0647848534332951
$data[20,9] = 0
$data[20,10] = 0.5236523606862136
$data[20,11] = 0.4423393052639213
$data[20,12] = 0
$data[20,13] = 2.854863957184705
$data[21,0] = 1.669394380660833
$data[21,1] = 0.3320439500865007
$data[21,2] = 0
$data[21,3] = 0.2211693817008644
$data[21,4] = 1.800027743670668
$data[21,5] = 0.002436904505353941
$data[21,6] = 0
$data[21,7] = 0.7514293301963413
$data[21,8] = 0.06374543107812514
$data[21,9] = 0
$data[21,10] = 0.5200723777308127
$data[21,11] = 0.432875672458465
$data[21,12] = 0
$data[21,13] = 2.861399602438723
$data[22,0] = 1.475938723551792
$data[22,1] = 0.3089969186212045
$data[22,2] = 0
$data[22,3] = 0.2229741923827753
$data[22,4] = 1.804679466846068
$data[22,5] = 0.002440940135898423
$data[22,6] = 0
$data[22,7] = 0.7675639320525924
$data[22,8] = 0.05980344470265919
$data[22,9] = 0
$data[22,10] = 0.506771276131488
$data[22,11] = 0.3971111762140538
$data[22,12] = 0
$data[22,13] = 2.889890973343739
$data[23,0] = 1.267069733252924
$data[23,1] = 0.2839509752973015
$data[23,2] = 0
$data[23,3] = 0.2251550558986093
$data[23,4] = 1.813655624748918
$data[23,5] = 0.002445621027913569
$data[23,6] = 0
$data[23,7] = 0.7869170167233293
$data[23,8] = 0.05554157966538753
$data[23,9] = 0
$data[23,10] = 0.4929809615317282
$data[23,11] = 0.3587120385928486
$data[23,12] = 0
$data[23,13] = 2.9285833456033

$ws.Range("B2:O25").Value = $data
